$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q2").Value = 257
$ws.Range("R2").Value = 339.9528082890812
$ws.Range("S2").Value = 2.297667413202893
$ws.Range("T2").Value = 12.5

$ws.Range("Q3").Value = 257
$ws.Range("R3").Value = 339.9528082890812
$ws.Range("S3").Value = 4.595334826405786
$ws.Range("T3").Value = 25

$ws.Range("Q4").Value = 257
$ws.Range("R4").Value = 339.9528082890812
$ws.Range("S4").Value = 38.59022489036856
$ws.Range("T4").Value = 209.9424

$ws.Range("Q5").Value = 257
$ws.Range("R5").Value = 339.9528082890812
$ws.Range("S5").Value = 113.0177740864556
$ws.Range("T5").Value = 614.8505949829331

$ws.Range("Q6").Value = 257
$ws.Range("R6").Value = 339.9528082890812
$ws.Range("S6").Value = 113.0177740864556
$ws.Range("T6").Value = 614.8505949829331

$ws.Range("Q7").Value = 257
$ws.Range("R7").Value = 339.9528082890812
$ws.Range("S7").Value = 113.0177740864556
$ws.Range("T7").Value = 614.8505949829331

$ws.Range("Q8").Value = 255.3500061035156
$ws.Range("R8").Value = 337.7702399670201
$ws.Range("S8").Value = 119.9385280045054
$ws.Range("T8").Value = 652.5015724388176

$ws.Range("Q9").Value = 255.3500061035156
$ws.Range("R9").Value = 337.7702399670201
$ws.Range("S9").Value = 359.2085170977269
$ws.Range("T9").Value = 1954.202091181894

$ws.Range("Q10").Value = 255.3500061035156
$ws.Range("R10").Value = 337.7702399670201
$ws.Range("S10").Value = 112.3193522233961
$ws.Range("T10").Value = 611.0509705298559

$ws.Range("Q11").Value = 255.3500061035156
$ws.Range("R11").Value = 337.7702399670201
$ws.Range("S11").Value = 111.4727771366061
$ws.Range("T11").Value = 606.4453480955267

$ws.Range("Q12").Value = 255.3500061035156
$ws.Range("R12").Value = 337.7702399670201
$ws.Range("S12").Value = 112.3193522233961
$ws.Range("T12").Value = 611.0509705298559

$ws.Range("Q13").Value = 257
$ws.Range("R13").Value = 339.9528082890812
$ws.Range("S13").Value = 108.784898652506
$ws.Range("T13").Value = 591.8224828112876

$ws.Range("Q14").Value = 257
$ws.Range("R14").Value = 339.9528082890812
$ws.Range("S14").Value = 108.784898652506
$ws.Range("T14").Value = 591.8224828112876

$ws.Range("Q15").Value = 257
$ws.Range("R15").Value = 339.9528082890812
$ws.Range("S15").Value = 108.784898652506
$ws.Range("T15").Value = 591.8224828112876

$ws.Range("Q16").Value = 250.6499938964844
$ws.Range("R16").Value = 331.5531880262682
$ws.Range("S16").Value = 110.3298956023555
$ws.Range("T16").Value = 600.227729698694

$ws.Range("Q17").Value = 250.6499938964844
$ws.Range("R17").Value = 331.5531880262682
$ws.Range("S17").Value = 352.5968688702955
$ws.Range("T17").Value = 1918.232741410907

$ws.Range("Q18").Value = 257
$ws.Range("R18").Value = 339.9528082890812
$ws.Range("S18").Value = 98.89055232564868
$ws.Range("T18").Value = 537.9942706100665

$ws.Range("Q19").Value = 257
$ws.Range("R19").Value = 339.9528082890812
$ws.Range("S19").Value = 119.3670872373801
$ws.Range("T19").Value = 649.3927632404012

$ws.Range("Q20").Value = 257
$ws.Range("R20").Value = 339.9528082890812
$ws.Range("S20").Value = 104
$ws.Range("T20").Value = 565.7911987304688

$ws.Range("Q21").Value = 257
$ws.Range("R21").Value = 339.9528082890812
$ws.Range("S21").Value = 104
$ws.Range("T21").Value = 565.7911987304688

$ws.Range("Q22").Value = 257
$ws.Range("R22").Value = 339.9528082890812
$ws.Range("S22").Value = 1.838133930562314
$ws.Range("T22").Value = 10

$ws.Range("Q23").Value = 257
$ws.Range("R23").Value = 339.9528082890812
$ws.Range("S23").Value = 1.571236883844666
$ws.Range("T23").Value = 8.548

$ws.Range("Q24").Value = 257
$ws.Range("R24").Value = 339.9528082890812
$ws.Range("S24").Value = 119.3670872373801
$ws.Range("T24").Value = 649.3927632404012

$ws.Range("Q25").Value = 257
$ws.Range("R25").Value = 339.9528082890812
$ws.Range("S25").Value = 113.0177740864556
$ws.Range("T25").Value = 614.8505949829331

$ws.Range("Q26").Value = 257
$ws.Range("R26").Value = 339.9528082890812
$ws.Range("S26").Value = 129.9492758222542
$ws.Range("T26").Value = 706.9630436695147

$ws.Range("Q27").Value = 257
$ws.Range("R27").Value = 339.9528082890812
$ws.Range("S27").Value = 112.1711989996657
$ws.Range("T27").Value = 610.244972548604

$ws.Range("Q28").Value = 255.3500061035156
$ws.Range("R28").Value = 337.7702399670201
$ws.Range("S28").Value = 359.2085170977269
$ws.Range("T28").Value = 1954.202091181894

$ws.Range("Q29").Value = 255.3500061035156
$ws.Range("R29").Value = 337.7702399670201
$ws.Range("S29").Value = 359.2085170977269
$ws.Range("T29").Value = 1954.202091181894

$ws.Range("S30").Value = 3611.665917223789
$ws.Range("T30").Value = 19648.54604538486
